$wb = $excel.ActiveWorkbook

# The localization status moved on from "Ready for handoff" to "In Translation".
# That text shows up in the Overview sheet (columns E/F hold the per-language
# status) and on each language sheet (column C, "Status").
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value2 = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value2 = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value2 = "In Translation"

# Narrow the "Status" columns to match the new (shorter) text's autofit width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
